$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 2 (done first): paragraph "After the designing of the bio
# page, I started working on the form page." + bookmark(_GoBack) + " "
# Merge back into a single run (with trailing space) and drop the
# old _GoBack bookmark from here - it relocates to the other edit
# below, exactly like Word moves its hidden "last edit" bookmark.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$afterPara = $d.Paragraphs.Item(9)
$afterRange = $afterPara.Range
$afterRange.End = $afterRange.End - 1
$afterRange.Delete()
$afterRange.InsertAfter("After the designing of the bio page, I started working on the form page. ")

# ------------------------------------------------------------------
# Change 1: paragraph "Fig:The basic woreframing of the form page."
# Split into "Fig" + "(3)" + bookmark(_GoBack) + ":The basic ..."
# ------------------------------------------------------------------
$figPara = $d.Paragraphs.Item(4)
$figStart = $figPara.Range.Start

# Insert "(3)" right after "Fig" (3 characters in)
$insertPoint = $d.Range($figStart + 3, $figStart + 3)
$insertPoint.InsertBefore("(3)")

# Force a run break between "Fig" and "(3)" (otherwise identical
# formatting would keep them coalesced into a single run) by toggling
# the underline off/on over the newly inserted "(3)" text.
$threeRange = $d.Range($figStart + 3, $figStart + 6)
$threeRange.Font.Underline = 0
$threeRange.Font.Underline = 1

# Insert the (re-appearing) _GoBack bookmark between "(3)" and the
# remaining ":The basic woreframing of the form page." text.
$bmRange = $d.Range($figStart + 6, $figStart + 6)
$d.Bookmarks.Add("_GoBack", $bmRange)
